$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from the existing "sum" header (G1) onto the
# new "Save" header cell (H1), then set its value.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add the new data value under the new header.
$ws.Range("H2").Value = 1
